$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each row that changed,
# matching the refreshed values from the GitHub Actions data pull.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.315.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.856.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.47%  "

$ws.Range("E4").Value = "  -0.82%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("E6").Value = "  -0.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4610"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.88%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07322"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8827"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.47%  "

$ws.Range("E11").Value = "  -0.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07796"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.872.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.394"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.546"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008921"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.338.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.122"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.081.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.917"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.067"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.125"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7650"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.001"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.172"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.39%  "

$ws.Range("E35").Value = "  +0.94%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.623"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01963"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.076"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05224"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.942"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.036"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5153"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1640"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.369"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4836"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.0000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.653"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06219"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.01%  "
